# ----------------------------------------------------------------------------
# Applies the "Fixed update to excel issue" commit:
#   1. Renames the "Requested quantity" header on the "Weekly Quantity" and
#      "Monthly Trend" sheets to "Weekly_PO_Qty" / "Monthly_PO_Qty".
#   2. Appends a brand-new "PO Forecast" worksheet (ds / PO_Forecast /
#      yhat_lower / yhat_upper) populated with the forecast series.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers -----------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet as the last tab ----------------------
# Duplicate an existing sheet (rather than Worksheets.Add()) so sheet-level
# properties (outline/page-setup defaults, margins, etc.) match the rest of
# the workbook, then wipe its contents before writing the forecast data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeekly.Copy([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.Clear()

# --- Column headers -----------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the bold/centered/thin-bordered header formatting already used on
# the other sheets instead of re-declaring fonts/borders by hand.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Reuse the existing date-time number format for column A.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Forecast data (rows 2-60) ----------------------------------------------
$wsForecast.Cells.Item(2,1).Value = 44948.99999999999; $wsForecast.Cells.Item(2,2).Value = 27; $wsForecast.Cells.Item(2,3).Value = -49.89490641956051; $wsForecast.Cells.Item(2,4).Value = 101.909528447338
$wsForecast.Cells.Item(3,1).Value = 44962.99999999999; $wsForecast.Cells.Item(3,2).Value = 28; $wsForecast.Cells.Item(3,3).Value = -51.31606627111479; $wsForecast.Cells.Item(3,4).Value = 112.8398467253717
$wsForecast.Cells.Item(4,1).Value = 44983.99999999999; $wsForecast.Cells.Item(4,2).Value = 30; $wsForecast.Cells.Item(4,3).Value = -50.03271220138149; $wsForecast.Cells.Item(4,4).Value = 109.6625579416308
$wsForecast.Cells.Item(5,1).Value = 44997.99999999999; $wsForecast.Cells.Item(5,2).Value = 32; $wsForecast.Cells.Item(5,3).Value = -42.54114083705946; $wsForecast.Cells.Item(5,4).Value = 110.7098854824107
$wsForecast.Cells.Item(6,1).Value = 45004.99999999999; $wsForecast.Cells.Item(6,2).Value = 32; $wsForecast.Cells.Item(6,3).Value = -51.38810454647714; $wsForecast.Cells.Item(6,4).Value = 107.7133664759144
$wsForecast.Cells.Item(7,1).Value = 45011.99999999999; $wsForecast.Cells.Item(7,2).Value = 33; $wsForecast.Cells.Item(7,3).Value = -45.03474075204286; $wsForecast.Cells.Item(7,4).Value = 107.0779524316076
$wsForecast.Cells.Item(8,1).Value = 45032.99999999999; $wsForecast.Cells.Item(8,2).Value = 35; $wsForecast.Cells.Item(8,3).Value = -42.61386050813437; $wsForecast.Cells.Item(8,4).Value = 109.2448052717468
$wsForecast.Cells.Item(9,1).Value = 45039.99999999999; $wsForecast.Cells.Item(9,2).Value = 36; $wsForecast.Cells.Item(9,3).Value = -41.24288548441663; $wsForecast.Cells.Item(9,4).Value = 110.9238315628198
$wsForecast.Cells.Item(10,1).Value = 45046.99999999999; $wsForecast.Cells.Item(10,2).Value = 36; $wsForecast.Cells.Item(10,3).Value = -39.04778056735849; $wsForecast.Cells.Item(10,4).Value = 110.8171020317638
$wsForecast.Cells.Item(11,1).Value = 45053.99999999999; $wsForecast.Cells.Item(11,2).Value = 37; $wsForecast.Cells.Item(11,3).Value = -43.25611666063916; $wsForecast.Cells.Item(11,4).Value = 113.7620481658415
$wsForecast.Cells.Item(12,1).Value = 45060.99999999999; $wsForecast.Cells.Item(12,2).Value = 38; $wsForecast.Cells.Item(12,3).Value = -41.2749439761702; $wsForecast.Cells.Item(12,4).Value = 116.5972261268139
$wsForecast.Cells.Item(13,1).Value = 45074.99999999999; $wsForecast.Cells.Item(13,2).Value = 39; $wsForecast.Cells.Item(13,3).Value = -29.24159212466549; $wsForecast.Cells.Item(13,4).Value = 116.5466428230553
$wsForecast.Cells.Item(14,1).Value = 45081.99999999999; $wsForecast.Cells.Item(14,2).Value = 40; $wsForecast.Cells.Item(14,3).Value = -37.85625813746935; $wsForecast.Cells.Item(14,4).Value = 116.6545243590602
$wsForecast.Cells.Item(15,1).Value = 45088.99999999999; $wsForecast.Cells.Item(15,2).Value = 40; $wsForecast.Cells.Item(15,3).Value = -36.69024125719639; $wsForecast.Cells.Item(15,4).Value = 112.827659620333
$wsForecast.Cells.Item(16,1).Value = 45095.99999999999; $wsForecast.Cells.Item(16,2).Value = 41; $wsForecast.Cells.Item(16,3).Value = -31.21918941550126; $wsForecast.Cells.Item(16,4).Value = 112.0631773910582
$wsForecast.Cells.Item(17,1).Value = 45102.99999999999; $wsForecast.Cells.Item(17,2).Value = 42; $wsForecast.Cells.Item(17,3).Value = -30.64664986305803; $wsForecast.Cells.Item(17,4).Value = 115.9067315714838
$wsForecast.Cells.Item(18,1).Value = 45109.99999999999; $wsForecast.Cells.Item(18,2).Value = 42; $wsForecast.Cells.Item(18,3).Value = -39.4751995546424; $wsForecast.Cells.Item(18,4).Value = 126.2061795680067
$wsForecast.Cells.Item(19,1).Value = 45116.99999999999; $wsForecast.Cells.Item(19,2).Value = 43; $wsForecast.Cells.Item(19,3).Value = -31.17846389244074; $wsForecast.Cells.Item(19,4).Value = 120.4132195169913
$wsForecast.Cells.Item(20,1).Value = 45123.99999999999; $wsForecast.Cells.Item(20,2).Value = 44; $wsForecast.Cells.Item(20,3).Value = -34.46682367933708; $wsForecast.Cells.Item(20,4).Value = 116.8136801992978
$wsForecast.Cells.Item(21,1).Value = 45130.99999999999; $wsForecast.Cells.Item(21,2).Value = 44; $wsForecast.Cells.Item(21,3).Value = -42.5757516763147; $wsForecast.Cells.Item(21,4).Value = 127.7953892929935
$wsForecast.Cells.Item(22,1).Value = 45137.99999999999; $wsForecast.Cells.Item(22,2).Value = 45; $wsForecast.Cells.Item(22,3).Value = -35.5452020368699; $wsForecast.Cells.Item(22,4).Value = 126.7711089423088
$wsForecast.Cells.Item(23,1).Value = 45144.99999999999; $wsForecast.Cells.Item(23,2).Value = 46; $wsForecast.Cells.Item(23,3).Value = -34.61601508029364; $wsForecast.Cells.Item(23,4).Value = 125.7203030825644
$wsForecast.Cells.Item(24,1).Value = 45151.99999999999; $wsForecast.Cells.Item(24,2).Value = 46; $wsForecast.Cells.Item(24,3).Value = -29.7035559832001; $wsForecast.Cells.Item(24,4).Value = 119.2032929041421
$wsForecast.Cells.Item(25,1).Value = 45158.99999999999; $wsForecast.Cells.Item(25,2).Value = 47; $wsForecast.Cells.Item(25,3).Value = -33.32011467807712; $wsForecast.Cells.Item(25,4).Value = 129.5203385018096
$wsForecast.Cells.Item(26,1).Value = 45165.99999999999; $wsForecast.Cells.Item(26,2).Value = 48; $wsForecast.Cells.Item(26,3).Value = -34.99403417672711; $wsForecast.Cells.Item(26,4).Value = 126.4162272224346
$wsForecast.Cells.Item(27,1).Value = 45179.99999999999; $wsForecast.Cells.Item(27,2).Value = 49; $wsForecast.Cells.Item(27,3).Value = -26.37043692461325; $wsForecast.Cells.Item(27,4).Value = 127.7366223483436
$wsForecast.Cells.Item(28,1).Value = 45193.99999999999; $wsForecast.Cells.Item(28,2).Value = 50; $wsForecast.Cells.Item(28,3).Value = -21.06901476261864; $wsForecast.Cells.Item(28,4).Value = 130.3851240583466
$wsForecast.Cells.Item(29,1).Value = 45200.99999999999; $wsForecast.Cells.Item(29,2).Value = 51; $wsForecast.Cells.Item(29,3).Value = -29.20111171345718; $wsForecast.Cells.Item(29,4).Value = 126.9459960643389
$wsForecast.Cells.Item(30,1).Value = 45207.99999999999; $wsForecast.Cells.Item(30,2).Value = 52; $wsForecast.Cells.Item(30,3).Value = -26.78233349041281; $wsForecast.Cells.Item(30,4).Value = 124.4373469815517
$wsForecast.Cells.Item(31,1).Value = 45214.99999999999; $wsForecast.Cells.Item(31,2).Value = 52; $wsForecast.Cells.Item(31,3).Value = -25.60434128498783; $wsForecast.Cells.Item(31,4).Value = 132.4107615851609
$wsForecast.Cells.Item(32,1).Value = 45221.99999999999; $wsForecast.Cells.Item(32,2).Value = 53; $wsForecast.Cells.Item(32,3).Value = -24.20865596292498; $wsForecast.Cells.Item(32,4).Value = 131.4833424285025
$wsForecast.Cells.Item(33,1).Value = 45235.99999999999; $wsForecast.Cells.Item(33,2).Value = 54; $wsForecast.Cells.Item(33,3).Value = -21.75845090227321; $wsForecast.Cells.Item(33,4).Value = 131.3759150665537
$wsForecast.Cells.Item(34,1).Value = 45242.99999999999; $wsForecast.Cells.Item(34,2).Value = 55; $wsForecast.Cells.Item(34,3).Value = -25.11695178559205; $wsForecast.Cells.Item(34,4).Value = 130.1560740185786
$wsForecast.Cells.Item(35,1).Value = 45249.99999999999; $wsForecast.Cells.Item(35,2).Value = 56; $wsForecast.Cells.Item(35,3).Value = -26.57669815423872; $wsForecast.Cells.Item(35,4).Value = 135.8970339215611
$wsForecast.Cells.Item(36,1).Value = 45256.99999999999; $wsForecast.Cells.Item(36,2).Value = 56; $wsForecast.Cells.Item(36,3).Value = -23.6167641217647; $wsForecast.Cells.Item(36,4).Value = 141.6583889214229
$wsForecast.Cells.Item(37,1).Value = 45263.99999999999; $wsForecast.Cells.Item(37,2).Value = 57; $wsForecast.Cells.Item(37,3).Value = -19.33126697019238; $wsForecast.Cells.Item(37,4).Value = 137.2183980846697
$wsForecast.Cells.Item(38,1).Value = 45270.99999999999; $wsForecast.Cells.Item(38,2).Value = 58; $wsForecast.Cells.Item(38,3).Value = -23.81338295140746; $wsForecast.Cells.Item(38,4).Value = 138.6922338855085
$wsForecast.Cells.Item(39,1).Value = 45277.99999999999; $wsForecast.Cells.Item(39,2).Value = 58; $wsForecast.Cells.Item(39,3).Value = -16.19077781882645; $wsForecast.Cells.Item(39,4).Value = 136.9651304145952
$wsForecast.Cells.Item(40,1).Value = 45298.99999999999; $wsForecast.Cells.Item(40,2).Value = 60; $wsForecast.Cells.Item(40,3).Value = -17.58087222865679; $wsForecast.Cells.Item(40,4).Value = 136.9793031250102
$wsForecast.Cells.Item(41,1).Value = 45396.99999999999; $wsForecast.Cells.Item(41,2).Value = 70; $wsForecast.Cells.Item(41,3).Value = -10.04312266915973; $wsForecast.Cells.Item(41,4).Value = 148.8420886642375
$wsForecast.Cells.Item(42,1).Value = 45417.99999999999; $wsForecast.Cells.Item(42,2).Value = 72; $wsForecast.Cells.Item(42,3).Value = -3.144210458688549; $wsForecast.Cells.Item(42,4).Value = 151.2718374433067
$wsForecast.Cells.Item(43,1).Value = 45424.99999999999; $wsForecast.Cells.Item(43,2).Value = 72; $wsForecast.Cells.Item(43,3).Value = -5.633035759282486; $wsForecast.Cells.Item(43,4).Value = 152.42369205045
$wsForecast.Cells.Item(44,1).Value = 45438.99999999999; $wsForecast.Cells.Item(44,2).Value = 74; $wsForecast.Cells.Item(44,3).Value = -6.047972382884622; $wsForecast.Cells.Item(44,4).Value = 147.7805581851242
$wsForecast.Cells.Item(45,1).Value = 45459.99999999999; $wsForecast.Cells.Item(45,2).Value = 76; $wsForecast.Cells.Item(45,3).Value = -0.7048697757631891; $wsForecast.Cells.Item(45,4).Value = 155.5872998983812
$wsForecast.Cells.Item(46,1).Value = 45466.99999999999; $wsForecast.Cells.Item(46,2).Value = 76; $wsForecast.Cells.Item(46,3).Value = -2.282084810069696; $wsForecast.Cells.Item(46,4).Value = 155.1876242031898
$wsForecast.Cells.Item(47,1).Value = 45487.99999999999; $wsForecast.Cells.Item(47,2).Value = 78; $wsForecast.Cells.Item(47,3).Value = 2.642002254232882; $wsForecast.Cells.Item(47,4).Value = 156.8084398052961
$wsForecast.Cells.Item(48,1).Value = 45515.99999999999; $wsForecast.Cells.Item(48,2).Value = 81; $wsForecast.Cells.Item(48,3).Value = 0.6315199050284386; $wsForecast.Cells.Item(48,4).Value = 158.377147679623
$wsForecast.Cells.Item(49,1).Value = 45529.99999999999; $wsForecast.Cells.Item(49,2).Value = 82; $wsForecast.Cells.Item(49,3).Value = 9.546789547348466; $wsForecast.Cells.Item(49,4).Value = 162.1547665090987
$wsForecast.Cells.Item(50,1).Value = 45543.99999999999; $wsForecast.Cells.Item(50,2).Value = 84; $wsForecast.Cells.Item(50,3).Value = 4.214253805073339; $wsForecast.Cells.Item(50,4).Value = 159.3003237015926
$wsForecast.Cells.Item(51,1).Value = 45606.99999999999; $wsForecast.Cells.Item(51,2).Value = 90; $wsForecast.Cells.Item(51,3).Value = 10.51676987940798; $wsForecast.Cells.Item(51,4).Value = 166.1628049194937
$wsForecast.Cells.Item(52,1).Value = 45613.99999999999; $wsForecast.Cells.Item(52,2).Value = 90; $wsForecast.Cells.Item(52,3).Value = 17.15441029620325; $wsForecast.Cells.Item(52,4).Value = 168.1251834002128
$wsForecast.Cells.Item(53,1).Value = 45620.99999999999; $wsForecast.Cells.Item(53,2).Value = 91; $wsForecast.Cells.Item(53,3).Value = 7.997882536495004; $wsForecast.Cells.Item(53,4).Value = 165.2555988200743
$wsForecast.Cells.Item(54,1).Value = 45627.99999999999; $wsForecast.Cells.Item(54,2).Value = 92; $wsForecast.Cells.Item(54,3).Value = 17.05040740663026; $wsForecast.Cells.Item(54,4).Value = 171.9695228707895
$wsForecast.Cells.Item(55,1).Value = 45634.99999999999; $wsForecast.Cells.Item(55,2).Value = 92; $wsForecast.Cells.Item(55,3).Value = 13.65978097776628; $wsForecast.Cells.Item(55,4).Value = 171.3790914716311
$wsForecast.Cells.Item(56,1).Value = 45641.99999999999; $wsForecast.Cells.Item(56,2).Value = 93; $wsForecast.Cells.Item(56,3).Value = 15.48231670878632; $wsForecast.Cells.Item(56,4).Value = 169.3324778429131
$wsForecast.Cells.Item(57,1).Value = 45648.99999999999; $wsForecast.Cells.Item(57,2).Value = 94; $wsForecast.Cells.Item(57,3).Value = 12.35805771607907; $wsForecast.Cells.Item(57,4).Value = 166.5294530541681
$wsForecast.Cells.Item(58,1).Value = 45655.99999999999; $wsForecast.Cells.Item(58,2).Value = 94; $wsForecast.Cells.Item(58,3).Value = 12.75950990735967; $wsForecast.Cells.Item(58,4).Value = 173.6608185928251
$wsForecast.Cells.Item(59,1).Value = 45662.99999999999; $wsForecast.Cells.Item(59,2).Value = 95; $wsForecast.Cells.Item(59,3).Value = 13.39032161284346; $wsForecast.Cells.Item(59,4).Value = 169.2833301749655
$wsForecast.Cells.Item(60,1).Value = 45669.99999999999; $wsForecast.Cells.Item(60,2).Value = 96; $wsForecast.Cells.Item(60,3).Value = 20.02274630533843; $wsForecast.Cells.Item(60,4).Value = 170.3240457420945

Write-Output "PO Forecast sheet added and headers renamed."
